# Automatic update: refreshed meteo.cat observation data + banner timestamps
# (commit: "Update automàtic: dades i banners [2026-02-05 08:41]")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks: columns M (13) and N (14) narrowed from 25 to 24 chars ---
# (23.15 is the COM ColumnWidth input that round-trips to an exact stored width of 24)
$ws.Columns.Item(13).ColumnWidth = 23.15
$ws.Columns.Item(14).ColumnWidth = 23.15

# --- Cell value refresh ---
# Most values are plain text already (station readings with units, dates, etc.)
# and can be written directly. A handful of humidity cells are bare "NN%"
# strings, which Excel would otherwise auto-convert to a numeric percentage -
# those are written with a leading apostrophe to force literal text, matching
# the workbook's original inline-string storage.

$ws.Range("E2").Value = '2026-02-05 08:39:39'
$ws.Range("E3").Value = '2026-02-05 08:39:42'
$ws.Range("G3").Value = '178 cm'
$ws.Range("H3").Value = '''63%'
$ws.Range("K3").Value = '-0.1 MJ/m2'
$ws.Range("L3").Value = '32.4 km/h - 219º 6:57 TU'
$ws.Range("M3").Value = '1.0 °C 6:46 TU'
$ws.Range("O3").Value = '-1.9 °C'
$ws.Range("E4").Value = '2026-02-05 08:39:44'
$ws.Range("H4").Value = '''74%'
$ws.Range("J4").Value = '993.5 hPa'
$ws.Range("K4").Value = '0.1 MJ/m2'
$ws.Range("M4").Value = '10.3 °C 8:29 TU'
$ws.Range("O4").Value = '8.2 °C'
$ws.Range("E5").Value = '2026-02-05 08:39:47'
$ws.Range("H5").Value = '''84%'
$ws.Range("J5").Value = '993.7 hPa'
$ws.Range("K5").Value = '0.2 MJ/m2'
$ws.Range("M5").Value = '6.9 °C 8:21 TU'
$ws.Range("N5").Value = '3.8 °C 7:00 TU'
$ws.Range("O5").Value = '5.1 °C'
$ws.Range("E6").Value = '2026-02-05 08:39:50'
$ws.Range("H6").Value = '''67%'
$ws.Range("J6").Value = '994.4 hPa'
$ws.Range("L6").Value = '14.4 km/h - 290º 7:42 TU'
$ws.Range("M6").Value = '12.6 °C 6:16 TU'
$ws.Range("O6").Value = '11.5 °C'
$ws.Range("E7").Value = '2026-02-05 08:39:52'
$ws.Range("H7").Value = '''72%'
$ws.Range("J7").Value = '994.3 hPa'
$ws.Range("L7").Value = '43.6 km/h - 321º 6:32 TU'
$ws.Range("M7").Value = '10.5 °C 6:52 TU'
$ws.Range("O7").Value = '9.6 °C'
$ws.Range("E8").Value = '2026-02-05 08:39:55'
$ws.Range("M8").Value = '5.0 °C 7:59 TU'
$ws.Range("O8").Value = '3.4 °C'
$ws.Range("E9").Value = '2026-02-05 08:39:57'
$ws.Range("E10").Value = '2026-02-05 08:39:59'
$ws.Range("O10").Value = '2.6 °C'
$ws.Range("E11").Value = '2026-02-05 08:40:02'
$ws.Range("J11").Value = '998.7 hPa'
$ws.Range("K11").Value = '0.2 MJ/m2'
$ws.Range("M11").Value = '-0.2 °C 8:29 TU'
$ws.Range("O11").Value = '-1.5 °C'
$ws.Range("E12").Value = '2026-02-05 08:40:05'
$ws.Range("H12").Value = '''92%'
$ws.Range("L12").Value = '21.2 km/h - 321º 7:08 TU'
$ws.Range("M12").Value = '8.3 °C 7:18 TU'
$ws.Range("O12").Value = '6.9 °C'
$ws.Range("E13").Value = '2026-02-05 08:40:08'
$ws.Range("O13").Value = '3.5 °C'
$ws.Range("E14").Value = '2026-02-05 08:40:11'
$ws.Range("G14").Value = '68 cm'
$ws.Range("H14").Value = '''75%'
$ws.Range("I14").Value = '0.0 mm'
$ws.Range("K14").Value = '0.0 MJ/m2'
$ws.Range("L14").Value = '50.0 km/h - 231º 0:25 TU'
$ws.Range("M14").Value = '-3.1 °C 0:10 TU'
$ws.Range("N14").Value = '-3.3 °C 0:27 TU'
$ws.Range("O14").Value = '-3.2 °C'
$ws.Range("E15").Value = '2026-02-05 08:40:13'
$ws.Range("E16").Value = '2026-02-05 08:40:16'
$ws.Range("E17").Value = '2026-02-05 08:40:19'
$ws.Range("J17").Value = '998.0 hPa'
$ws.Range("L17").Value = '12.2 km/h - 94º 6:59 TU'
$ws.Range("M17").Value = '1.0 °C 7:59 TU'
$ws.Range("O17").Value = '0.0 °C'
$ws.Range("E18").Value = '2026-02-05 08:40:22'
$ws.Range("N18").Value = '-5.6 °C 3:33 TU'
$ws.Range("O18").Value = '-5.1 °C'
$ws.Range("E19").Value = '2026-02-05 08:40:24'
$ws.Range("I19").Value = '0.6 mm'
$ws.Range("J19").Value = '995.6 hPa'
$ws.Range("L19").Value = '13.7 km/h - 86º 5:44 TU'
$ws.Range("M19").Value = '6.2 °C 6:13 TU'
$ws.Range("O19").Value = '4.5 °C'
$ws.Range("E20").Value = '2026-02-05 08:40:27'
$ws.Range("G20").Value = '111 cm'
$ws.Range("H20").Value = '''63%'
$ws.Range("K20").Value = '-0.1 MJ/m2'
$ws.Range("L20").Value = '29.5 km/h - 241º 5:48 TU'
$ws.Range("M20").Value = '0.0 °C 5:54 TU'
$ws.Range("N20").Value = '-4.6 °C 6:07 TU'
$ws.Range("O20").Value = '-2.1 °C'
$ws.Range("E21").Value = '2026-02-05 08:40:29'
$ws.Range("E22").Value = '2026-02-05 08:40:32'
$ws.Range("E23").Value = '2026-02-05 08:40:35'
$ws.Range("E24").Value = '2026-02-05 08:40:38'
$ws.Range("H24").Value = '''85%'
$ws.Range("J24").Value = '993.1 hPa'
$ws.Range("K24").Value = '0.1 MJ/m2'
$ws.Range("L24").Value = '22.0 km/h - 203º 7:52 TU'
$ws.Range("O24").Value = '8.3 °C'
$ws.Range("E25").Value = '2026-02-05 08:40:40'
$ws.Range("E26").Value = '2026-02-05 08:40:43'
$ws.Range("G26").Value = '110 cm'
$ws.Range("H26").Value = '''67%'
$ws.Range("K26").Value = '0.1 MJ/m2'
$ws.Range("O26").Value = '-3.2 °C'
$ws.Range("E27").Value = '2026-02-05 08:40:46'
$ws.Range("I27").Value = '0.1 mm'
$ws.Range("J27").Value = '993.7 hPa'
$ws.Range("K27").Value = '0.1 MJ/m2'
$ws.Range("N27").Value = '3.1 °C 5:49 TU'
$ws.Range("O27").Value = '4.1 °C'
$ws.Range("E28").Value = '2026-02-05 08:40:49'
$ws.Range("J28").Value = '998.1 hPa'
$ws.Range("L28").Value = '20.5 km/h - 256º 6:43 TU'
$ws.Range("M28").Value = '-0.4 °C 6:24 TU'
$ws.Range("O28").Value = '-1.4 °C'
$ws.Range("E29").Value = '2026-02-05 08:40:52'
$ws.Range("H29").Value = '''92%'
$ws.Range("K29").Value = '0.1 MJ/m2'
$ws.Range("L29").Value = '16.9 km/h - 306º 7:09 TU'
$ws.Range("M29").Value = '6.7 °C 8:29 TU'
$ws.Range("O29").Value = '4.6 °C'
$ws.Range("E30").Value = '2026-02-05 08:40:54'
$ws.Range("H30").Value = '''58%'
$ws.Range("K30").Value = '-0.1 MJ/m2'
$ws.Range("N30").Value = '-5.2 °C 3:40 TU'
$ws.Range("O30").Value = '-3.0 °C'
$ws.Range("E31").Value = '2026-02-05 08:40:57'
$ws.Range("G31").Value = '0 cm'
$ws.Range("E32").Value = '2026-02-05 08:41:00'
$ws.Range("H32").Value = '''86%'
$ws.Range("I32").Value = '0.4 mm'
$ws.Range("J32").Value = '994.3 hPa'
$ws.Range("O32").Value = '8.9 °C'
$ws.Range("E33").Value = '2026-02-05 08:41:02'
$ws.Range("N33").Value = '3.0 °C 5:32 TU'
$ws.Range("O33").Value = '3.7 °C'
$ws.Range("E34").Value = '2026-02-05 08:41:05'
$ws.Range("E35").Value = '2026-02-05 08:41:08'
$ws.Range("E36").Value = '2026-02-05 08:41:10'
$ws.Range("I36").Value = '1.0 mm'
$ws.Range("J36").Value = '995.6 hPa'
$ws.Range("M36").Value = '6.4 °C 7:51 TU'
$ws.Range("O36").Value = '5.4 °C'
